# petty-cashBook-2021.xlsx — 15-Jan-2021 midday update
# Adds the 13-Jan, 14-Jan and 15-Jan (serials 44209/44210/44211) petty-cash
# transactions to the "Buku KAS HARIAN"-style ledger on Sheet1, continuing
# on from the existing 12-Jan (44208) entries, and moves the frozen-pane
# selection down to where entry was left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# ---- 12-Jan-2021 (A14, serial 44208) continued ----
$ws.Range("B14").Value = "Wages Expense"
$ws.Range("D14").Formula = "=45000+210000"

$ws.Range("B15").Value = "TRANSFER BCA"
$ws.Range("D15").Formula = "=855000+3210000+445000+120000+2223000+45705000+100000+2200000+155000+3207000+336000"

$ws.Range("B16").Value = "A/R"
$ws.Range("C16").Formula = "=45705000+25212000"

$ws.Range("B17").Value = "FREIGHT OUT"
$ws.Range("D17").Formula = "=204000"

$ws.Range("B18").Value = "SALES - cash/retail"
$ws.Range("C18").Formula = "=12692475+22178525-25212000"

$ws.Range("B19").Value = "SELISIH - kurang"
$ws.Range("D19").Value = 81000

$ws.Range("B20").Value = "SETOR KE BANK"
$ws.Range("D20").Value = 21000000

# ---- 13-Jan-2021 (A21, serial 44209) ----
$ws.Range("A21").Value = 44209
$ws.Range("A21").NumberFormat = $ws.Range("A14").NumberFormat

$ws.Range("B21").Value = "Wages Expense"
$ws.Range("D21").Formula = "=45000+255000"

$ws.Range("B22").Value = "TRANSFER BCA"
$ws.Range("D22").Formula = "=36000000+12000000+285000+320000+757000+9775000+3850000+3247000+34200000+937500+715000"

$ws.Range("B23").Value = "FREIGHT OUT"
$ws.Range("D23").Formula = "=56000+56000+28000"

$ws.Range("B24").Value = "A/R"
$ws.Range("C24").Formula = "=12000000+52700000+34200000+78056500"

$ws.Range("B25").Value = "Undangan - RITA"
$ws.Range("D25").Formula = "=5000000"

$ws.Range("B26").Value = "BELI kresek"
$ws.Range("D26").Formula = "=50000"

$ws.Range("B27").Value = "PLN - Astar 214"
$ws.Range("D27").Formula = "=103000"

$ws.Range("B28").Value = "SALES - cash/retail"
$ws.Range("C28").Formula = "=7636975+76158025-78056500"

$ws.Range("B29").Value = "SELISIH - lebih"
$ws.Range("C29").Value = 48000

$ws.Range("B30").Value = "SETOR KE BANK"
$ws.Range("D30").Formula = "=75000000"

# ---- 14-Jan-2021 (A31, serial 44210) ----
$ws.Range("A31").Value = 44210
$ws.Range("A31").NumberFormat = $ws.Range("A14").NumberFormat

$ws.Range("B31").Value = "Wages Expense"
$ws.Range("D31").Formula = "=45000+270000"

$ws.Range("B32").Value = "A/R"
$ws.Range("C32").Formula = "=20000000+8197500+2000000+39483000"

$ws.Range("B33").Value = "TRANSFER BCA"
$ws.Range("D33").Formula = "=1225000+7721000+839000+100000"

$ws.Range("B34").Value = "PLN - Astar 165"
$ws.Range("D34").Formula = "=815000"

$ws.Range("B35").Value = "Telpon - 5224823"
$ws.Range("D35").Value = 252000

$ws.Range("B36").Value = "prive - andreas"
$ws.Range("D36").Formula = "=5000000"

$ws.Range("B37").Value = "SALES - cash/retail"
$ws.Range("C37").Formula = "=63341525-15251525-39483000"

$ws.Range("B38").Value = "SELISIH - kurang"
$ws.Range("D38").Value = 60000

$ws.Range("B39").Value = "SETOR KE BANK"
$ws.Range("D39").Value = 62000000

# ---- 15-Jan-2021 (A40, serial 44211) ----
$ws.Range("A40").Value = 44211
$ws.Range("A40").NumberFormat = $ws.Range("A14").NumberFormat

$ws.Range("B40").Value = "Wages Expense"
$ws.Range("D40").Formula = "=45000"

$ws.Range("B41").Value = "A/R"
$ws.Range("C41").Formula = "=37292000+2308000"

$ws.Range("B42").Value = "TRANSFER BCA"
$ws.Range("D42").Formula = "=39600000+4800000"

# ---- move the frozen-pane view down to where today's entry stopped ----
$win = $excel.ActiveWindow
$win.ScrollRow = 39
$ws.Range("E60").Select()
